# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.749.82'
$ws.Range('E2').Value = '  +0.36%  '

# Row 3
$ws.Range('D3').Value = '1.650.49'
$ws.Range('E3').Value = '  +0.85%  '

# Row 4
$ws.Range('E4').Value = '  -0.05%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.86'
$ws.Range('E5').Value = '  +1.22%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.506'
$ws.Range('E6').Value = '  +1.36%  '

# Row 7
$ws.Range('E7').Value = '  -0.01%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.255'
$ws.Range('E8').Value = '  -0.13%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0631'
$ws.Range('E9').Value = '  +1.20%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.46'
$ws.Range('E10').Value = '  +1.70%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0846'
$ws.Range('E11').Value = '  +0.64%  '

# Row 12
$ws.Range('D12').Value = '1.880.31'
$ws.Range('E12').Value = '  +0.78%  '

# Row 13
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.23'
$ws.Range('E13').Value = '  +3.24%  '

# Row 14
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.649.52'
$ws.Range('E14').Value = '  +0.73%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.537'
$ws.Range('E15').Value = '  +1.69%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.79'
$ws.Range('E16').Value = '  +5.52%  '

# Row 17
$ws.Range('D17').Value = '26.782.61'
$ws.Range('E17').Value = '  +0.44%  '

# Row 18
$ws.Range('D18').Value = '0.0₃0758'
$ws.Range('E18').Value = '  +1.80%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '222.15'
$ws.Range('E19').Value = '  +1.72%  '

# Row 20
$ws.Range('E20').Value = '  +0.04%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.42'
$ws.Range('E21').Value = '  +2.49%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.39'
$ws.Range('E22').Value = '  +2.65%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.61'
$ws.Range('E23').Value = '  +1.09%  '

# Row 24
$ws.Range('E24').Value = '  +12.00%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.45'
$ws.Range('E25').Value = '  -1.11%  '

# Row 26
$ws.Range('E26').Value = '  +0.00%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.122'
$ws.Range('E27').Value = '  +0.44%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.14'
$ws.Range('E28').Value = '  +4.07%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.99'
$ws.Range('E29').Value = '  +3.62%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0525'
$ws.Range('E30').Value = '  +1.72%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.17'
$ws.Range('E31').Value = '  +0.55%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.45'
$ws.Range('E32').Value = '  +4.48%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.08'
$ws.Range('E33').Value = '  +4.63%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.58'
$ws.Range('E34').Value = '  +4.21%  '

# Row 35
$ws.Range('D35').Value = '1.296.59'
$ws.Range('E35').Value = '  +8.66%  '

# Row 36
$ws.Range('E36').Value = '  +6.01%  '

# Row 37
$ws.Range('E37').Value = '  +0.88%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.835'
$ws.Range('E38').Value = '  +3.22%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.530'
$ws.Range('E39').Value = '  +4.68%  '

# Row 40
$ws.Range('E40').Value = '  +0.05%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.815'
$ws.Range('E41').Value = '  +2.55%  '

# Row 42
$ws.Range('E42').Value = '  -2.91%  '

# Row 43
$ws.Range('E43').Value = '  +0.67%  '

# Row 44
$ws.Range('D44').Value = '1.791.66'
$ws.Range('E44').Value = '  +1.11%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '93.96'
$ws.Range('E45').Value = '  +1.89%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '60.37'
$ws.Range('E46').Value = '  +10.06%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.63'
$ws.Range('E47').Value = '  +5.84%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0519'
$ws.Range('E48').Value = '  +1.19%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.80'
$ws.Range('E49').Value = '  +1.54%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0984'
$ws.Range('E50').Value = '  +3.91%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.407'
$ws.Range('E51').Value = '  -0.76%  '
